$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 5 ("state-ref") for the new "basestate" entry.
# This pushes all following rows down by one (old row 5 -> row 6, etc.)
$ws.Rows.Item(5).Insert()

# Fill in the newly inserted row 5, matching the style/section of rows 3-4 ("state" group)
$ws.Cells.Item(5, 2).Value = "basestate"

# In the row that now holds "command" (row 8 after the insert), add "exit" in column E
$ws.Cells.Item(8, 5).Value = "exit"

# Update the active selection to match the authored file
$ws.Range("E9").Select()
